# REPORTGEN-397 : update extension name for quality standards
#
# 1. Split " id from quality model (eg page" so that "eg" sits in its own
#    run (mirrors Word's automatic spell-check run-splitting around the
#    abbreviation "eg").
# 2. Replace the extension name "Standard Quality Rules" with
#    "Quality Standards Support" inside the "**" footnote paragraph, and
#    drop Word's "last edit" (_GoBack) bookmark right after the newly
#    typed text - exactly where Word leaves it after an interactive edit.
#    Adding a new _GoBack bookmark automatically retires the previous one
#    (Word only ever keeps a single _GoBack bookmark in a document), which
#    takes care of removing the old bookmark that used to sit around the
#    chart picture further down the document.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: "... (eg page 4)" -> split "eg" into its own run
# ---------------------------------------------------------------------
$egFind = $d.Content
$egFind.Find.Execute("(eg page", $true, $false, $false, $false, $false, `
                      $true, 1, $false, "", 0) | Out-Null

$egStart = $egFind.Start + 1          # skip the leading "("
$egEnd = $egStart + 2                 # "eg" is two characters
$egRange = $d.Range($egStart, $egEnd)

if ($egRange.Text -eq "eg") {
    # Round-tripping a character-formatting property forces the run to be
    # split out from its identically-formatted neighbours without leaving
    # any visible/semantic formatting difference behind.
    $egRange.Font.Bold = $true
    $egRange.Font.Bold = $false
    Write-Output "Change 1 applied: 'eg' isolated into its own run"
} else {
    Write-Output "Change 1 SKIPPED: could not locate 'eg' in '(eg page'"
}

# ---------------------------------------------------------------------
# Change 2: extension name rename + _GoBack bookmark relocation
# ---------------------------------------------------------------------
$introFind = $d.Content
$introFind.Find.Execute("** The selection of metrics", $true, $false, `
                         $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$phraseRange = $d.Range($introFind.End, $d.Content.End)
$phraseRange.Find.Execute("Standard Quality Rules", $true, $false, $false, `
                           $false, $false, $true, 1, $false, "", 0) | Out-Null

if ($phraseRange.Text -eq "Standard Quality Rules") {
    $newTextStart = $phraseRange.Start
    $phraseRange.Text = "Quality Standards Support"
    $newTextEnd = $newTextStart + "Quality Standards Support".Length

    # Split the freshly typed text into its own run, same trick as above.
    $newTextRange = $d.Range($newTextStart, $newTextEnd)
    $newTextRange.Font.Bold = $true
    $newTextRange.Font.Bold = $false

    # Drop the _GoBack bookmark right after the inserted text - this is
    # where Word leaves it after the last edit, and adding it here removes
    # it from its old location automatically.
    $bookmarkPoint = $d.Range($newTextEnd, $newTextEnd)
    $d.Bookmarks.Add("_GoBack", $bookmarkPoint) | Out-Null
    Write-Output "Change 2 applied: extension name updated and _GoBack bookmark relocated"
} else {
    Write-Output "Change 2 SKIPPED: could not locate 'Standard Quality Rules'"
}
